$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 120 - this shifts the existing rows 120:169 down to 121:170,
# carrying their values/styles along (matches the diff: new row170 == old row169, etc.)
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with this week's record. Columns A,B,C,E,F,G,I,Q,R
# keep the same constant values used throughout this data block.
$ws.Cells.Item(120, 1).Value = 4
$ws.Cells.Item(120, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(120, 3).Value = "Los Lagos"
$ws.Cells.Item(120, 4).Value = 45146
$ws.Cells.Item(120, 5).Value = 10
$ws.Cells.Item(120, 6).Value = 100112022
$ws.Cells.Item(120, 7).Value = "Arveja Verde"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 35
$ws.Cells.Item(120, 11).Value = 39000
$ws.Cells.Item(120, 12).Value = 39000
$ws.Cells.Item(120, 13).Value = 39000
$ws.Cells.Item(120, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(120, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(120, 16).Value = 1560
$ws.Cells.Item(120, 17).Value = 25
$ws.Cells.Item(120, 18).Value = "Hortaliza"
